$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.1574606916619176
$wsLinear.Range("B3").Value = 0.1090491499087974
$wsLinear.Range("B4").Value = 51.08948806223462
$wsLinear.Range("B5").Value = "[1.0, 0.14504445694985793, 0.025472600447849082, 0.08613002374663972, 0.08407760926014012, 0.0185559692342061, 0.08377570569365057, 0.27270368940927353, 0.10016201233836572, -0.014444711809234873, 0.09706648172344313, 0.07321634639226914, -0.024113806525724014, 0.09179266925151944, 0.2667902887524999, 0.06590626011605515, -0.012790364047846687, 0.11000259250597443, 0.09547334351363555, -0.025292533265979234]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -0.01642966077857282
$wsNonLinear.Range("B5").Value = 0.2565507903589642
$wsNonLinear.Range("B6").Value = 48.50647356568494
$wsNonLinear.Range("B7").Value = -0.4378328735474312
$wsNonLinear.Range("B8").Value = 0.6760901503594756
$wsNonLinear.Range("B9").Value = 53.58126562322729
$wsNonLinear.Range("B10").Value = "[1.0, 0.14418144470406055, 0.02484336659429966, 0.08627028257247382, 0.08450115467944377, 0.018162349315222648, 0.08239256422519736, 0.2720375516792358, 0.0995402779258603, -0.014797074853573554, 0.09722934224706116, 0.07324057602205596, -0.024199142171614865, 0.09216411055487878, 0.2657882729432112, 0.06559405304487399, -0.013058476709091611, 0.109844474652561, 0.09540798768701135, -0.025674479940512984]"

Write-Host "Edits applied"
